# Apply card allotment / market display copy-count changes.
$wb = $excel.ActiveWorkbook

# StarterDeck: ST03 "Weak Rocket" quantity 2 -> 1
$wsStarter = $wb.Worksheets.Item("StarterDeck")
$wsStarter.Range("F4").Value = 1

# RocketMarket: copies column (G) for rows 2, 4, 8 go from 2 -> 3
$wsRocket = $wb.Worksheets.Item("RocketMarket")
$wsRocket.Range("G2").Value = 3
$wsRocket.Range("G4").Value = 3
$wsRocket.Range("G8").Value = 3

# ShieldMarket: copies column (G) for rows 3, 4 go from 2 -> 3
$wsShield = $wb.Worksheets.Item("ShieldMarket")
$wsShield.Range("G3").Value = 3
$wsShield.Range("G4").Value = 3
